# Finalization of the Reports through GA4 and Chartbeat
#
# The sheet currently has: Site | Total Users (GA4) | Views (GA4) | Views per
# user | Engagement rate | Average engagement time  (cols A-F)
#
# Three new columns are inserted after "Site" (B) for the Chartbeat /
# Adsense numbers, pushing the existing GA4 columns to E-I:
#   A Site
#   B Adsense Revenue            (new, blank data cells)
#   C Total users (chartbeat)    (new)
#   D Views (chartbeat)          (new)
#   E Total Users (GA4)          (was B)
#   F Views (GA4)                (was C)
#   G Views per user             (was D)
#   H Engagement rate            (was E)
#   I Average engagement time    (was F)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns in front of the old column B (Total Users (GA4)).
# This shifts B:F -> E:I, the row/cell data moves with it, and the sheet
# dimension grows from A1:F5 to A1:I5 automatically.
$ws.Range("B1:D1").EntireColumn.Insert()

# New header row (row 2) cells for the inserted columns.
$ws.Range("B2").Value = "Adsense Revenue"
$ws.Range("C2").Value = "Total users (chartbeat)"
$ws.Range("D2").Value = "Views (chartbeat)"

# Chartbeat users/views per site (Adsense Revenue column is left blank).
$ws.Range("C3").Value = "153.619"
$ws.Range("D3").Value = "88.845"

$ws.Range("C4").Value = "136.241"
$ws.Range("D4").Value = "71.161"

$ws.Range("C5").Value = "6.581"
$ws.Range("D5").Value = "4.745"

# Column widths. A keeps its width; the two new "chartbeat" name columns
# (B, C) get the wide 25.83-char width, D keeps the 15.83-char width that
# used to belong to old column B, and E:I re-use the widths that used to
# belong to old columns C:F (shifted along with the data).
$ws.Columns.Item(2).ColumnWidth = 25   # B - Adsense Revenue
$ws.Columns.Item(3).ColumnWidth = 25   # C - Total users (chartbeat)
$ws.Columns.Item(4).ColumnWidth = 15   # D - Views (chartbeat)
$ws.Columns.Item(5).ColumnWidth = 12   # E - Total Users (GA4)
$ws.Columns.Item(6).ColumnWidth = 15   # F - Views (GA4)
$ws.Columns.Item(7).ColumnWidth = 15   # G - Views per user
$ws.Columns.Item(8).ColumnWidth = 20   # H - Engagement rate
$ws.Columns.Item(9).ColumnWidth = 35   # I - Average engagement time
